$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.380.63'
$ws.Range("E2").Value = '  +2.30%  '
$ws.Range("D3").Value = '2.092.75'
$ws.Range("E3").Value = '  -0.03%  '
$ws.Range("E4").Value = '  -0.91%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '342.97'
$ws.Range("D5").Style = "Normal"
$ws.Range("E6").Value = '  -0.81%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5244'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.71%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4420'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.93%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '54.54'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.65%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.09319'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.68%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.168'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.03%  '
$ws.Range("E12").Value = '  +0.17%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '8.583'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.84%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '2.117.85'
$ws.Range("E14").Value = '  +0.74%  '
$ws.Range("B15").Value = 'Polkadot'
$ws.Range("C15").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.898'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.45%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '101.27'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.10%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001156'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.61%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.001'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.90%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '21.10'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.69%  '
$ws.Range("E20").Value = '  -0.10%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.322'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.19%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.0000'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.88%  '
$ws.Range("D23").Value = '30.391.79'
$ws.Range("E23").Value = '  +2.21%  '
$ws.Range("E24").Value = '  +0.11%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.311'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.42%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '21.79'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.39%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '162.87'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.78%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.500'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.55%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '133.13'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.27%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.138'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.79%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.669'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.56%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1045'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.56%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.869'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +10.12%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.246'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.50%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.862'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.85%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '10.13'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.83%  '
$ws.Range("E37").Value = '  +2.76%  '
$ws.Range("E38").Value = '  +2.30%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.6983'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.75%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '12.54'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.92%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.342'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.96%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.2214'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.25%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.6809'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.12%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.40'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.49%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.346'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.60%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.000'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.77%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.376'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +18.42%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.631'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.18%  '
$ws.Range("E49").Value = '  -1.69%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.215'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +8.80%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.215'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.23%  '
